$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new geolocation ping rows captured a bit later the same evening
# (next test run, ~17:05-17:08 UTC) are appended after the current last row.
$newRows = @(
    @{
        A = "202.173.124.126"
        B = 28.3621617
        C = 77.2827806
        D = 15.01099967956543
        E = "Mozilla/5.0 (Linux; Android 10; K) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Mobile Safari/537.36"
        F = "Linux armv81"
        G = "2025-06-25T17:05:51.318Z"
    },
    @{
        A = "202.173.124.126"
        B = 28.475392
        C = 77.0670592
        D = 616570.7228211587
        E = "Mozilla/5.0 (Windows NT 10.0; Win64; x64) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Safari/537.36"
        F = "Win32"
        G = "2025-06-25T17:05:53.723Z"
    },
    @{
        A = "202.173.124.126"
        B = 28.3621537
        C = 77.2828149
        D = 12.38599967956543
        E = "Mozilla/5.0 (Linux; Android 10; K) AppleWebKit/537.36 (KHTML, like Gecko) Chrome/137.0.0.0 Mobile Safari/537.36"
        F = "Linux armv81"
        G = "2025-06-25T17:07:58.748Z"
    }
)

# Find the current bottom of the data table so this keeps working even if
# the sheet already has a different number of rows.
$lastRow = $ws.UsedRange.Rows.Count
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
}

$endRow = $startRow + $newRows.Count - 1

# Keep the "number stored as text" error-checking suppression in sync with
# the (now larger) data range, matching the sheet's existing convention.
try {
    $fullRange = $ws.Range("A1:G$endRow")
    $fullRange.Errors.Item(3).Ignore = $true
} catch {
}
